$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "DauerStunden" header was actually minutes all along -> rename to DauerMinuten
$ws.Range("F1").Value = "DauerMinuten"

# Append the new arrival/departure scan rows (26-46) that were recorded since
# the export was last generated.
$ws.Range("A26").Value = "Eli"
$ws.Range("B26").Value = "Enders"
$ws.Range("C26").Value = "4a"
$ws.Range("D26").Value = "12.11.2023 14:00"
$ws.Range("E26").Value = "12.11.2023 14:00"
$ws.Range("F26").Value = 0

$ws.Range("A27").Value = "Max"
$ws.Range("B27").Value = "Schmitz"
$ws.Range("C27").Value = "4a"
$ws.Range("D27").Value = "12.11.2023 14:03"
$ws.Range("E27").Value = "12.11.2023 14:11"
$ws.Range("F27").Value = 8

$ws.Range("A28").Value = "Stephan"
$ws.Range("B28").Value = "Fuchs"
$ws.Range("C28").Value = "3C"
$ws.Range("D28").Value = "12.11.2023 14:04"
$ws.Range("E28").Value = "12.11.2023 14:11"
$ws.Range("F28").Value = 7

$ws.Range("A29").Value = "Detlef"
$ws.Range("B29").Value = "Soost"
$ws.Range("C29").Value = "1a"
$ws.Range("D29").Value = "12.11.2023 14:08"
$ws.Range("E29").Value = "12.11.2023 14:11"
$ws.Range("F29").Value = 3

$ws.Range("A30").Value = "Eli"
$ws.Range("B30").Value = "Enders"
$ws.Range("C30").Value = "4a"
$ws.Range("D30").Value = "13.11.2023 21:49"
$ws.Range("E30").Value = "13.11.2023 21:49"
$ws.Range("F30").Value = 0

$ws.Range("A31").Value = "Max"
$ws.Range("B31").Value = "Schmitz"
$ws.Range("C31").Value = "4a"
$ws.Range("D31").Value = "13.11.2023 21:51"
$ws.Range("E31").Value = "13.11.2023 21:51"
$ws.Range("F31").Value = 0

$ws.Range("A32").Value = "Stephan"
$ws.Range("B32").Value = "Fuchs"
$ws.Range("C32").Value = "3C"
$ws.Range("D32").Value = "13.11.2023 21:58"
$ws.Range("E32").Value = "13.11.2023 21:58"
$ws.Range("F32").Value = 0

$ws.Range("A33").Value = "Stephan"
$ws.Range("B33").Value = "Fuchs"
$ws.Range("C33").Value = "3C"
$ws.Range("D33").Value = "14.11.2023 17:29"
$ws.Range("E33").Value = "14.11.2023 17:29"
$ws.Range("F33").Value = 0

$ws.Range("A34").Value = "Max"
$ws.Range("B34").Value = "Schmitz"
$ws.Range("C34").Value = "4a"
$ws.Range("D34").Value = "19.11.2023 18:54"
$ws.Range("E34").Value = "19.11.2023 19:34"
$ws.Range("F34").Value = 40

$ws.Range("A35").Value = "Eli"
$ws.Range("B35").Value = "Enders"
$ws.Range("C35").Value = "4a"
$ws.Range("D35").Value = "19.11.2023 19:35"
$ws.Range("E35").Value = "19.11.2023 19:35"
$ws.Range("F35").Value = 0

$ws.Range("A36").Value = "Eli"
$ws.Range("B36").Value = "Enders"
$ws.Range("C36").Value = "4a"
$ws.Range("D36").Value = "22.11.2023 10:32"
$ws.Range("E36").Value = "22.11.2023 10:32"
$ws.Range("F36").Value = 0

$ws.Range("A37").Value = "Detlef"
$ws.Range("B37").Value = "Soost"
$ws.Range("C37").Value = "1a"
$ws.Range("D37").Value = "22.11.2023 13:53"
$ws.Range("E37").Value = "22.11.2023 13:53"
$ws.Range("F37").Value = 0

$ws.Range("A38").Value = "Stephan"
$ws.Range("B38").Value = "Fuchs"
$ws.Range("C38").Value = "3C"
$ws.Range("D38").Value = "26.11.2023 16:02"
$ws.Range("E38").Value = "26.11.2023 16:02"
$ws.Range("F38").Value = 0

$ws.Range("A39").Value = "Stephan"
$ws.Range("B39").Value = "Fuchs"
$ws.Range("C39").Value = "3C"
$ws.Range("D39").Value = "07.12.2023 20:40"
$ws.Range("E39").Value = "07.12.2023 20:41"
$ws.Range("F39").Value = 1

$ws.Range("A40").Value = "Eli"
$ws.Range("B40").Value = "Enders"
$ws.Range("C40").Value = "4a"
$ws.Range("D40").Value = "07.12.2023 20:40"
$ws.Range("E40").Value = "07.12.2023 20:41"
$ws.Range("F40").Value = 1

$ws.Range("A41").Value = "Stephan"
$ws.Range("B41").Value = "Fuchs"
$ws.Range("C41").Value = "3C"
$ws.Range("D41").Value = "17.12.2023 22:06"
$ws.Range("E41").Value = "17.12.2023 22:06"
$ws.Range("F41").Value = 0

$ws.Range("A42").Value = "Eli"
$ws.Range("B42").Value = "Enders"
$ws.Range("C42").Value = "4a"
$ws.Range("D42").Value = "17.12.2023 22:06"
$ws.Range("E42").Value = "17.12.2023 22:06"
$ws.Range("F42").Value = 0

$ws.Range("A43").Value = "Detlef"
$ws.Range("B43").Value = "Soost"
$ws.Range("C43").Value = "1a"
$ws.Range("D43").Value = "17.12.2023 22:06"
$ws.Range("E43").Value = "17.12.2023 22:06"
$ws.Range("F43").Value = 0

$ws.Range("A44").Value = "Stephan"
$ws.Range("B44").Value = "Fuchs"
$ws.Range("C44").Value = "3C"
$ws.Range("D44").Value = "18.12.2023 08:00"
$ws.Range("E44").Value = "18.12.2023 08:08"
$ws.Range("F44").Value = 8

$ws.Range("A45").Value = "Max"
$ws.Range("B45").Value = "Schmitz"
$ws.Range("C45").Value = "4a"
$ws.Range("D45").Value = "18.12.2023 08:00"
$ws.Range("E45").Value = "18.12.2023 08:01"
$ws.Range("F45").Value = 1

$ws.Range("A46").Value = "Detlef"
$ws.Range("B46").Value = "Soost"
$ws.Range("C46").Value = "1a"
$ws.Range("D46").Value = "18.12.2023 08:00"
$ws.Range("E46").Value = "18.12.2023 08:01"
$ws.Range("F46").Value = 1
